# DLC chapter overview - add "The Walking Dead" chapter row, clear the
# highlighted-row fill on row 3, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the green highlight fill that was previously on row 3 -------
# (A3:I3 currently carry a solid FF92D050 fill; the edit clears it back to
# "no fill" while leaving the rest of the formatting untouched.)
$ws.Range("A3:I3").Style = "Normal"

# --- 2. Append the new DLC row (row 47) -------------------------------------
$ws.Range("A47").Value = 59
$ws.Range("B47").Value = "The Walking Dead"
$ws.Range("C47").Value = "29.07.2025"
$ws.Range("D47").Value = 2
$ws.Range("E47").Formula = '=CHOOSE(D47, "Chapter DLC", "Half-Chapter DLC", "Clothing Pack DLC", "Original Soundtrack DLC", "Character Pack DLC", "Other", "Retracted", "Chapter Pack DLC")'
$ws.Range("F47").Value = "9.1.0"
$ws.Range("G47").Value = "Lasagna"
$ws.Range("H47").Value = "Lasagna"
$ws.Range("I47").Value = "Rick Grimes; Michonne Grimes"

# --- 3. Update the view: scroll back to the top and select C3 --------------
$ws.Range("A1").Select() | Out-Null
$ws.Range("C3").Select() | Out-Null
